# Auto-generated edit script applying scheduled market-price/profit refresh
# to the Coeurl_Profits workbook across all 8 crafting-class sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 1122311.4  # H6 1010090.2 -> 1122311.4
$ws.Cells.Item(6, 9).Value = 1122311.4  # I6 1010090.2 -> 1122311.4
$ws.Cells.Item(6, 11).Value = 3366934.2  # K6 3030270.6 -> 3366934.2
$ws.Cells.Item(6, 13).Value = -3366822.2  # M6 -3030158.6 -> -3366822.2
$ws.Cells.Item(9, 8).Value = 15228  # H9 13063.286 -> 15228
$ws.Cells.Item(9, 9).Value = 15228  # I9 13063.286 -> 15228
$ws.Cells.Item(9, 11).Value = 15228  # K9 13063.286 -> 15228
$ws.Cells.Item(9, 13).Value = -15059  # M9 -12894.286 -> -15059
$ws.Cells.Item(12, 8).Value = 244  # H12 268.33334 -> 244
$ws.Cells.Item(12, 9).Value = 124.5  # I12 133.33333 -> 124.5
$ws.Cells.Item(12, 11).Value = 124.5  # K12 133.33333 -> 124.5
$ws.Cells.Item(12, 13).Value = 45.5  # M12 36.66667000000001 -> 45.5
$ws.Cells.Item(15, 8).Value = 213.125  # H15 168.14285 -> 213.125
$ws.Cells.Item(15, 9).Value = 213.125  # I15 168.14285 -> 213.125
$ws.Cells.Item(15, 11).Value = 639.375  # K15 504.42855 -> 639.375
$ws.Cells.Item(15, 13).Value = -470.375  # M15 -335.42855 -> -470.375
$ws.Cells.Item(17, 8).Value = 2317177  # H17 6174639.5 -> 2317177
$ws.Cells.Item(17, 10).Value = 2317177  # J17 6174639.5 -> 2317177
$ws.Cells.Item(17, 12).Value = 6951531  # L17 18523918.5 -> 6951531
$ws.Cells.Item(17, 14).Value = -6951867  # N17 -18524254.5 -> -6951867
$ws.Cells.Item(29, 8).Value = 680.3333  # H29 4090.5 -> 680.3333
$ws.Cells.Item(29, 9).Value = 395.75  # I29 396.5 -> 395.75
$ws.Cells.Item(29, 10).Value = 1249.5  # J29 5937.5 -> 1249.5
$ws.Cells.Item(29, 11).Value = 1187.25  # K29 1189.5 -> 1187.25
$ws.Cells.Item(29, 12).Value = 3748.5  # L29 17812.5 -> 3748.5
$ws.Cells.Item(29, 13).Value = -906.25  # M29 -908.5 -> -906.25
$ws.Cells.Item(29, 14).Value = -4310.5  # N29 -18374.5 -> -4310.5
$ws.Cells.Item(38, 8).Value = 1274.8  # H38 725.3333 -> 1274.8
$ws.Cells.Item(38, 9).Value = 1274.8  # I38 725.3333 -> 1274.8
$ws.Cells.Item(38, 11).Value = 3824.4  # K38 2175.9999 -> 3824.4
$ws.Cells.Item(38, 13).Value = -3452.4  # M38 -1803.9999 -> -3452.4
$ws.Cells.Item(40, 8).Value = 4395.5557  # H40 4106 -> 4395.5557
$ws.Cells.Item(40, 9).Value = 4770.625  # I40 4407.222 -> 4770.625
$ws.Cells.Item(40, 11).Value = 4770.625  # K40 4407.222 -> 4770.625
$ws.Cells.Item(40, 13).Value = -4595.625  # M40 -4232.222 -> -4595.625
$ws.Cells.Item(42, 8).Value = 209.3  # H42 235.90909 -> 209.3
$ws.Cells.Item(42, 9).Value = 42.2  # I42 57.4 -> 42.2
$ws.Cells.Item(42, 10).Value = 376.4  # J42 384.66666 -> 376.4
$ws.Cells.Item(42, 11).Value = 126.6  # K42 172.2 -> 126.6
$ws.Cells.Item(42, 12).Value = 1129.2  # L42 1153.99998 -> 1129.2
$ws.Cells.Item(42, 13).Value = 103.4  # M42 57.80000000000001 -> 103.4
$ws.Cells.Item(42, 14).Value = -1589.2  # N42 -1613.99998 -> -1589.2
$ws.Cells.Item(51, 8).Value = 2638.5186  # H51 2636.1155 -> 2638.5186
$ws.Cells.Item(51, 10).Value = 2915.5557  # J51 2893.2222 -> 2915.5557
$ws.Cells.Item(51, 12).Value = 2915.5557  # L51 2893.2222 -> 2915.5557
$ws.Cells.Item(51, 14).Value = -3883.5557  # N51 -3861.2222 -> -3883.5557
$ws.Cells.Item(52, 8).Value = 0  # H52 977.5 -> 0
$ws.Cells.Item(52, 9).Value = 0  # I52 934.5 -> 0
$ws.Cells.Item(52, 10).Value = 0  # J52 999 -> 0
$ws.Cells.Item(52, 11).Value = 0  # K52 2803.5 -> 0
$ws.Cells.Item(52, 12).Value = 0  # L52 2997 -> 0
$ws.Cells.Item(52, 13).Value = $null  # clear M52 (was -2643.5)
$ws.Cells.Item(52, 14).Value = $null  # clear N52 (was -3317)
$ws.Cells.Item(53, 8).Value = 228.14285  # H53 227.42857 -> 228.14285
$ws.Cells.Item(53, 9).Value = 147  # I53 140.66667 -> 147
$ws.Cells.Item(53, 10).Value = 336.33334  # J53 383.6 -> 336.33334
$ws.Cells.Item(53, 11).Value = 147  # K53 140.66667 -> 147
$ws.Cells.Item(53, 12).Value = 336.33334  # L53 383.6 -> 336.33334
$ws.Cells.Item(53, 13).Value = 490  # M53 496.33333 -> 490
$ws.Cells.Item(53, 14).Value = -1610.33334  # N53 -1657.6 -> -1610.33334
$ws.Cells.Item(63, 8).Value = 70000  # H63 54999.5 -> 70000
$ws.Cells.Item(63, 9).Value = 0  # I63 29999 -> 0
$ws.Cells.Item(63, 10).Value = 70000  # J63 80000 -> 70000
$ws.Cells.Item(63, 11).Value = 0  # K63 29999 -> 0
$ws.Cells.Item(63, 12).Value = 70000  # L63 80000 -> 70000
$ws.Cells.Item(63, 13).Value = $null  # clear M63 (was -29375)
$ws.Cells.Item(63, 14).Value = -71248  # N63 -81248 -> -71248
$ws.Cells.Item(66, 8).Value = 70000  # H66 54999.5 -> 70000
$ws.Cells.Item(66, 9).Value = 0  # I66 29999 -> 0
$ws.Cells.Item(66, 10).Value = 70000  # J66 80000 -> 70000
$ws.Cells.Item(66, 11).Value = 0  # K66 89997 -> 0
$ws.Cells.Item(66, 12).Value = 210000  # L66 240000 -> 210000
$ws.Cells.Item(66, 13).Value = $null  # clear M66 (was -86877)
$ws.Cells.Item(66, 14).Value = -216240  # N66 -246240 -> -216240
$ws.Cells.Item(70, 8).Value = 75176.42999999999  # H70 80890.08 -> 75176.42999999999
$ws.Cells.Item(70, 9).Value = 1480.4  # I70 2300 -> 1480.4
$ws.Cells.Item(70, 10).Value = 116118.664  # J70 95179.17999999999 -> 116118.664
$ws.Cells.Item(70, 11).Value = 4441.200000000001  # K70 6900 -> 4441.200000000001
$ws.Cells.Item(70, 12).Value = 348355.992  # L70 285537.54 -> 348355.992
$ws.Cells.Item(70, 13).Value = -4171.200000000001  # M70 -6630 -> -4171.200000000001
$ws.Cells.Item(70, 14).Value = -348895.992  # N70 -286077.54 -> -348895.992
$ws.Cells.Item(73, 8).Value = 75176.42999999999  # H73 80890.08 -> 75176.42999999999
$ws.Cells.Item(73, 9).Value = 1480.4  # I73 2300 -> 1480.4
$ws.Cells.Item(73, 10).Value = 116118.664  # J73 95179.17999999999 -> 116118.664
$ws.Cells.Item(73, 11).Value = 4441.200000000001  # K73 6900 -> 4441.200000000001
$ws.Cells.Item(73, 12).Value = 348355.992  # L73 285537.54 -> 348355.992
$ws.Cells.Item(73, 13).Value = -3505.200000000001  # M73 -5964 -> -3505.200000000001
$ws.Cells.Item(73, 14).Value = -350227.992  # N73 -287409.54 -> -350227.992
$ws.Cells.Item(80, 8).Value = 933.8333  # H80 1041.125 -> 933.8333
$ws.Cells.Item(80, 9).Value = 610.6  # I80 876.6667 -> 610.6
$ws.Cells.Item(80, 10).Value = 1058.1538  # J80 1079.0769 -> 1058.1538
$ws.Cells.Item(80, 11).Value = 1831.8  # K80 2630.0001 -> 1831.8
$ws.Cells.Item(80, 12).Value = 3174.4614  # L80 3237.2307 -> 3174.4614
$ws.Cells.Item(80, 13).Value = -833.8000000000002  # M80 -1632.0001 -> -833.8000000000002
$ws.Cells.Item(80, 14).Value = -5170.4614  # N80 -5233.2307 -> -5170.4614
$ws.Cells.Item(83, 8).Value = 933.8333  # H83 1041.125 -> 933.8333
$ws.Cells.Item(83, 9).Value = 610.6  # I83 876.6667 -> 610.6
$ws.Cells.Item(83, 10).Value = 1058.1538  # J83 1079.0769 -> 1058.1538
$ws.Cells.Item(83, 11).Value = 5495.400000000001  # K83 7890.0003 -> 5495.400000000001
$ws.Cells.Item(83, 12).Value = 9523.3842  # L83 9711.6921 -> 9523.3842
$ws.Cells.Item(83, 13).Value = -503.4000000000005  # M83 -2898.0003 -> -503.4000000000005
$ws.Cells.Item(83, 14).Value = -19507.3842  # N83 -19695.6921 -> -19507.3842
$ws.Cells.Item(86, 8).Value = 5398  # H86 4684.2856 -> 5398
$ws.Cells.Item(86, 9).Value = 4331.6665  # I86 3759 -> 4331.6665
$ws.Cells.Item(86, 11).Value = 4331.6665  # K86 3759 -> 4331.6665
$ws.Cells.Item(86, 13).Value = -3208.6665  # M86 -2636 -> -3208.6665
$ws.Cells.Item(88, 8).Value = 1200.8  # H88 1216.2 -> 1200.8
$ws.Cells.Item(88, 9).Value = 0  # I88 777 -> 0
$ws.Cells.Item(88, 10).Value = 1200.8  # J88 1326 -> 1200.8
$ws.Cells.Item(88, 11).Value = 0  # K88 777 -> 0
$ws.Cells.Item(88, 12).Value = 1200.8  # L88 1326 -> 1200.8
$ws.Cells.Item(88, 13).Value = $null  # clear M88 (was -371)
$ws.Cells.Item(88, 14).Value = -2012.8  # N88 -2138 -> -2012.8
$ws.Cells.Item(89, 8).Value = 5398  # H89 4684.2856 -> 5398
$ws.Cells.Item(89, 9).Value = 4331.6665  # I89 3759 -> 4331.6665
$ws.Cells.Item(89, 11).Value = 21658.3325  # K89 18795 -> 21658.3325
$ws.Cells.Item(89, 13).Value = -16042.3325  # M89 -13179 -> -16042.3325
$ws.Cells.Item(91, 8).Value = 1200.8  # H91 1216.2 -> 1200.8
$ws.Cells.Item(91, 9).Value = 0  # I91 777 -> 0
$ws.Cells.Item(91, 10).Value = 1200.8  # J91 1326 -> 1200.8
$ws.Cells.Item(91, 11).Value = 0  # K91 777 -> 0
$ws.Cells.Item(91, 12).Value = 1200.8  # L91 1326 -> 1200.8
$ws.Cells.Item(91, 13).Value = $null  # clear M91 (was 627)
$ws.Cells.Item(91, 14).Value = -4008.8  # N91 -4134 -> -4008.8
$ws.Cells.Item(112, 8).Value = 31397.111  # H112 30575.514 -> 31397.111
$ws.Cells.Item(112, 9).Value = 1316.5714  # I112 1276.75 -> 1316.5714
$ws.Cells.Item(112, 11).Value = 3949.7142  # K112 3830.25 -> 3949.7142
$ws.Cells.Item(112, 13).Value = -2841.7142  # M112 -2722.25 -> -2841.7142
$ws.Cells.Item(113, 8).Value = 14774.723  # H113 16206.6875 -> 14774.723
$ws.Cells.Item(113, 9).Value = 12554.375  # I113 15632.833 -> 12554.375
$ws.Cells.Item(113, 11).Value = 12554.375  # K113 15632.833 -> 12554.375
$ws.Cells.Item(113, 13).Value = -9300.375  # M113 -12378.833 -> -9300.375
$ws.Cells.Item(132, 8).Value = 2722.2778  # H132 1994.7693 -> 2722.2778
$ws.Cells.Item(132, 9).Value = 2470.2942  # I132 1744.3334 -> 2470.2942
$ws.Cells.Item(132, 10).Value = 7006  # J132 5000 -> 7006
$ws.Cells.Item(132, 11).Value = 7410.882599999999  # K132 5233.0002 -> 7410.882599999999
$ws.Cells.Item(132, 12).Value = 21018  # L132 15000 -> 21018
$ws.Cells.Item(132, 13).Value = -4880.882599999999  # M132 -2703.0002 -> -4880.882599999999
$ws.Cells.Item(132, 14).Value = -26078  # N132 -20060 -> -26078
$ws.Cells.Item(137, 8).Value = 1951.2413  # H137 1859.1875 -> 1951.2413
$ws.Cells.Item(137, 9).Value = 1899.5333  # I137 1744.5 -> 1899.5333
$ws.Cells.Item(137, 11).Value = 5698.5999  # K137 5233.5 -> 5698.5999
$ws.Cells.Item(137, 13).Value = -3148.5999  # M137 -2683.5 -> -3148.5999
$ws.Cells.Item(138, 8).Value = 5321779  # H138 5378981.5 -> 5321779
$ws.Cells.Item(138, 9).Value = 1168.1482  # I138 1168.2222 -> 1168.1482
$ws.Cells.Item(138, 10).Value = 7465906  # J138 7578996.5 -> 7465906
$ws.Cells.Item(138, 11).Value = 3504.4446  # K138 3504.6666 -> 3504.4446
$ws.Cells.Item(138, 12).Value = 22397718  # L138 22736989.5 -> 22397718
$ws.Cells.Item(138, 13).Value = 1635.5554  # M138 1635.3334 -> 1635.5554
$ws.Cells.Item(138, 14).Value = -22407998  # N138 -22747269.5 -> -22407998

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(35, 8).Value = 2407.3333  # H35 3057.3333 -> 2407.3333
$ws.Cells.Item(35, 9).Value = 2407.3333  # I35 3057.3333 -> 2407.3333
$ws.Cells.Item(35, 11).Value = 2407.3333  # K35 3057.3333 -> 2407.3333
$ws.Cells.Item(35, 13).Value = -2001.3333  # M35 -2651.3333 -> -2001.3333
$ws.Cells.Item(61, 8).Value = 3659.182  # H61 3077.6428 -> 3659.182
$ws.Cells.Item(61, 9).Value = 2843.875  # I61 2326.0908 -> 2843.875
$ws.Cells.Item(61, 11).Value = 2843.875  # K61 2326.0908 -> 2843.875
$ws.Cells.Item(61, 13).Value = -2631.875  # M61 -2114.0908 -> -2631.875
$ws.Cells.Item(74, 8).Value = 7588.2188  # H74 7147.6763 -> 7588.2188
$ws.Cells.Item(74, 9).Value = 1127.25  # I74 1083.4073 -> 1127.25
$ws.Cells.Item(74, 10).Value = 26971.125  # J74 30538.428 -> 26971.125
$ws.Cells.Item(74, 11).Value = 1127.25  # K74 1083.4073 -> 1127.25
$ws.Cells.Item(74, 12).Value = 26971.125  # L74 30538.428 -> 26971.125
$ws.Cells.Item(74, 13).Value = -253.25  # M74 -209.4073000000001 -> -253.25
$ws.Cells.Item(74, 14).Value = -28719.125  # N74 -32286.428 -> -28719.125
$ws.Cells.Item(77, 8).Value = 7588.2188  # H77 7147.6763 -> 7588.2188
$ws.Cells.Item(77, 9).Value = 1127.25  # I77 1083.4073 -> 1127.25
$ws.Cells.Item(77, 10).Value = 26971.125  # J77 30538.428 -> 26971.125
$ws.Cells.Item(77, 11).Value = 5636.25  # K77 5417.0365 -> 5636.25
$ws.Cells.Item(77, 12).Value = 134855.625  # L77 152692.14 -> 134855.625
$ws.Cells.Item(77, 13).Value = -1268.25  # M77 -1049.0365 -> -1268.25
$ws.Cells.Item(77, 14).Value = -143591.625  # N77 -161428.14 -> -143591.625
$ws.Cells.Item(88, 8).Value = 1029.6666  # H88 1039.2778 -> 1029.6666
$ws.Cells.Item(88, 9).Value = 1054.7778  # I88 1178.2858 -> 1054.7778
$ws.Cells.Item(88, 10).Value = 992  # J88 950.8182 -> 992
$ws.Cells.Item(88, 11).Value = 1054.7778  # K88 1178.2858 -> 1054.7778
$ws.Cells.Item(88, 12).Value = 992  # L88 950.8182 -> 992
$ws.Cells.Item(88, 13).Value = -648.7778000000001  # M88 -772.2858000000001 -> -648.7778000000001
$ws.Cells.Item(88, 14).Value = -1804  # N88 -1762.8182 -> -1804
$ws.Cells.Item(91, 8).Value = 1029.6666  # H91 1039.2778 -> 1029.6666
$ws.Cells.Item(91, 9).Value = 1054.7778  # I91 1178.2858 -> 1054.7778
$ws.Cells.Item(91, 10).Value = 992  # J91 950.8182 -> 992
$ws.Cells.Item(91, 11).Value = 1054.7778  # K91 1178.2858 -> 1054.7778
$ws.Cells.Item(91, 12).Value = 992  # L91 950.8182 -> 992
$ws.Cells.Item(91, 13).Value = 349.2221999999999  # M91 225.7141999999999 -> 349.2221999999999
$ws.Cells.Item(91, 14).Value = -3800  # N91 -3758.8182 -> -3800
$ws.Cells.Item(97, 8).Value = 41864.6  # H97 49166.715 -> 41864.6
$ws.Cells.Item(97, 9).Value = 1753.5625  # I97 1939.8572 -> 1753.5625
$ws.Cells.Item(97, 10).Value = 113173.11  # J97 143620.42 -> 113173.11
$ws.Cells.Item(97, 11).Value = 1753.5625  # K97 1939.8572 -> 1753.5625
$ws.Cells.Item(97, 12).Value = 113173.11  # L97 143620.42 -> 113173.11
$ws.Cells.Item(97, 13).Value = -1257.5625  # M97 -1443.8572 -> -1257.5625
$ws.Cells.Item(97, 14).Value = -114165.11  # N97 -144612.42 -> -114165.11
$ws.Cells.Item(102, 8).Value = 2577.524  # H102 2763.238 -> 2577.524
$ws.Cells.Item(102, 9).Value = 2491.4666  # I102 2621.4666 -> 2491.4666
$ws.Cells.Item(102, 10).Value = 2792.6667  # J102 3117.6667 -> 2792.6667
$ws.Cells.Item(102, 11).Value = 2491.4666  # K102 2621.4666 -> 2491.4666
$ws.Cells.Item(102, 12).Value = 2792.6667  # L102 3117.6667 -> 2792.6667
$ws.Cells.Item(102, 13).Value = -869.4666000000002  # M102 -999.4666000000002 -> -869.4666000000002
$ws.Cells.Item(102, 14).Value = -6036.6667  # N102 -6361.6667 -> -6036.6667
$ws.Cells.Item(110, 8).Value = 9693.200000000001  # H110 9115.5 -> 9693.200000000001
$ws.Cells.Item(110, 9).Value = 11279.8  # I110 10295.272 -> 11279.8
$ws.Cells.Item(110, 11).Value = 11279.8  # K110 10295.272 -> 11279.8
$ws.Cells.Item(110, 13).Value = -9234.799999999999  # M110 -8250.272000000001 -> -9234.799999999999
$ws.Cells.Item(122, 8).Value = 1612.2106  # H122 1869.3572 -> 1612.2106
$ws.Cells.Item(122, 9).Value = 1334.5333  # I122 1555.7 -> 1334.5333
$ws.Cells.Item(122, 11).Value = 4003.5999  # K122 4667.1 -> 4003.5999
$ws.Cells.Item(122, 13).Value = -1553.5999  # M122 -2217.1 -> -1553.5999
$ws.Cells.Item(132, 8).Value = 3192.2222  # H132 3257.353 -> 3192.2222
$ws.Cells.Item(132, 9).Value = 3064.9395  # I132 3128.1614 -> 3064.9395
$ws.Cells.Item(132, 11).Value = 9194.818499999999  # K132 9384.484199999999 -> 9194.818499999999
$ws.Cells.Item(132, 13).Value = -6664.818499999999  # M132 -6854.484199999999 -> -6664.818499999999
$ws.Cells.Item(136, 8).Value = 3659.182  # H136 3077.6428 -> 3659.182
$ws.Cells.Item(136, 9).Value = 2843.875  # I136 2326.0908 -> 2843.875
$ws.Cells.Item(136, 11).Value = 8531.625  # K136 6978.2724 -> 8531.625
$ws.Cells.Item(136, 13).Value = -5981.625  # M136 -4428.2724 -> -5981.625

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(11, 8).Value = 215.75  # H11 252.375 -> 215.75
$ws.Cells.Item(11, 9).Value = 121  # I11 169.83333 -> 121
$ws.Cells.Item(11, 11).Value = 121  # K11 169.83333 -> 121
$ws.Cells.Item(11, 13).Value = 19  # M11 -29.83332999999999 -> 19
$ws.Cells.Item(134, 8).Value = 1817.575  # H134 1827.2208 -> 1817.575
$ws.Cells.Item(134, 9).Value = 1784.5065  # I134 1809.2933 -> 1784.5065
$ws.Cells.Item(134, 10).Value = 2666.3333  # J134 2499.5 -> 2666.3333
$ws.Cells.Item(134, 11).Value = 5353.5195  # K134 5427.8799 -> 5353.5195
$ws.Cells.Item(134, 12).Value = 7998.999899999999  # L134 7498.5 -> 7998.999899999999
$ws.Cells.Item(134, 13).Value = -2818.5195  # M134 -2892.8799 -> -2818.5195
$ws.Cells.Item(134, 14).Value = -13068.9999  # N134 -12568.5 -> -13068.9999

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(23, 8).Value = 36183.4  # H23 22990 -> 36183.4
$ws.Cells.Item(23, 10).Value = 44979.25  # J23 44980 -> 44979.25
$ws.Cells.Item(23, 12).Value = 44979.25  # L23 44980 -> 44979.25
$ws.Cells.Item(23, 14).Value = -45459.25  # N23 -45460 -> -45459.25
$ws.Cells.Item(27, 8).Value = 36183.4  # H27 22990 -> 36183.4
$ws.Cells.Item(27, 10).Value = 44979.25  # J27 44980 -> 44979.25
$ws.Cells.Item(27, 12).Value = 44979.25  # L27 44980 -> 44979.25
$ws.Cells.Item(27, 14).Value = -45363.25  # N27 -45364 -> -45363.25
$ws.Cells.Item(31, 8).Value = 65938.375  # H31 70310.92999999999 -> 65938.375
$ws.Cells.Item(31, 9).Value = 86395.414  # I31 94081.37 -> 86395.414
$ws.Cells.Item(31, 10).Value = 4567.25  # J31 4942.25 -> 4567.25
$ws.Cells.Item(31, 11).Value = 86395.414  # K31 94081.37 -> 86395.414
$ws.Cells.Item(31, 12).Value = 4567.25  # L31 4942.25 -> 4567.25
$ws.Cells.Item(31, 13).Value = -86100.414  # M31 -93786.37 -> -86100.414
$ws.Cells.Item(31, 14).Value = -5157.25  # N31 -5532.25 -> -5157.25
$ws.Cells.Item(34, 8).Value = 65938.375  # H34 70310.92999999999 -> 65938.375
$ws.Cells.Item(34, 9).Value = 86395.414  # I34 94081.37 -> 86395.414
$ws.Cells.Item(34, 10).Value = 4567.25  # J34 4942.25 -> 4567.25
$ws.Cells.Item(34, 11).Value = 86395.414  # K34 94081.37 -> 86395.414
$ws.Cells.Item(34, 12).Value = 4567.25  # L34 4942.25 -> 4567.25
$ws.Cells.Item(34, 13).Value = -86193.414  # M34 -93879.37 -> -86193.414
$ws.Cells.Item(34, 14).Value = -4971.25  # N34 -5346.25 -> -4971.25
$ws.Cells.Item(132, 8).Value = 4303.6875  # H132 4347.857 -> 4303.6875
$ws.Cells.Item(132, 9).Value = 4386  # I132 4541.222 -> 4386
$ws.Cells.Item(132, 10).Value = 4166.5  # J132 3999.8 -> 4166.5
$ws.Cells.Item(132, 11).Value = 13158  # K132 13623.666 -> 13158
$ws.Cells.Item(132, 12).Value = 12499.5  # L132 11999.4 -> 12499.5
$ws.Cells.Item(132, 13).Value = -10628  # M132 -11093.666 -> -10628
$ws.Cells.Item(132, 14).Value = -17559.5  # N132 -17059.4 -> -17559.5
$ws.Cells.Item(141, 8).Value = 235467  # H141 411715 -> 235467
$ws.Cells.Item(141, 10).Value = 235467  # J141 411715 -> 235467
$ws.Cells.Item(141, 12).Value = 235467  # L141 411715 -> 235467
$ws.Cells.Item(141, 14).Value = -245827  # N141 -422075 -> -245827

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(8, 8).Value = 212.5  # H8 0 -> 212.5
$ws.Cells.Item(8, 9).Value = 212.5  # I8 0 -> 212.5
$ws.Cells.Item(8, 11).Value = 637.5  # K8 0 -> 637.5
$ws.Cells.Item(8, 13).Value = -498.5  # M8 None -> -498.5
$ws.Cells.Item(49, 8).Value = 1099.1666  # H49 990.8333 -> 1099.1666
$ws.Cells.Item(49, 9).Value = 648.75  # I49 486.25 -> 648.75
$ws.Cells.Item(49, 11).Value = 1946.25  # K49 1458.75 -> 1946.25
$ws.Cells.Item(49, 13).Value = -1790.25  # M49 -1302.75 -> -1790.25
$ws.Cells.Item(92, 8).Value = 787.9048  # H92 770.13635 -> 787.9048
$ws.Cells.Item(92, 10).Value = 977.3570999999999  # J92 938.6667 -> 977.3570999999999
$ws.Cells.Item(92, 12).Value = 2932.0713  # L92 2816.0001 -> 2932.0713
$ws.Cells.Item(92, 14).Value = -5428.0713  # N92 -5312.0001 -> -5428.0713
$ws.Cells.Item(107, 8).Value = 518.8823  # H107 538.9375 -> 518.8823
$ws.Cells.Item(107, 9).Value = 431.25  # I107 452.45456 -> 431.25
$ws.Cells.Item(107, 11).Value = 1293.75  # K107 1357.36368 -> 1293.75
$ws.Cells.Item(107, 13).Value = 626.25  # M107 562.6363200000001 -> 626.25
$ws.Cells.Item(131, 8).Value = 21955.9  # H131 22849.896 -> 21955.9
$ws.Cells.Item(131, 9).Value = 333666.66  # I131 1000000 -> 333666.66
$ws.Cells.Item(131, 11).Value = 1000999.98  # K131 3000000 -> 1000999.98
$ws.Cells.Item(131, 13).Value = -995959.98  # M131 -2994960 -> -995959.98
$ws.Cells.Item(137, 8).Value = 3419.7273  # H137 3768.5557 -> 3419.7273
$ws.Cells.Item(137, 9).Value = 2043  # I137 2764.5 -> 2043
$ws.Cells.Item(137, 10).Value = 3936  # J137 4055.4285 -> 3936
$ws.Cells.Item(137, 11).Value = 6129  # K137 8293.5 -> 6129
$ws.Cells.Item(137, 12).Value = 11808  # L137 12166.2855 -> 11808
$ws.Cells.Item(137, 13).Value = -1029  # M137 -3193.5 -> -1029
$ws.Cells.Item(137, 14).Value = -22008  # N137 -22366.2855 -> -22008
$ws.Cells.Item(140, 8).Value = 3714.4546  # H140 4081.111 -> 3714.4546
$ws.Cells.Item(140, 9).Value = 2984.3333  # I140 3341.25 -> 2984.3333
$ws.Cells.Item(140, 10).Value = 7000  # J140 10000 -> 7000
$ws.Cells.Item(140, 11).Value = 8952.999899999999  # K140 10023.75 -> 8952.999899999999
$ws.Cells.Item(140, 12).Value = 21000  # L140 30000 -> 21000
$ws.Cells.Item(140, 13).Value = -3772.999899999999  # M140 -4843.75 -> -3772.999899999999
$ws.Cells.Item(140, 14).Value = -31360  # N140 -40360 -> -31360

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 7730.6  # H70 8034 -> 7730.6
$ws.Cells.Item(70, 10).Value = 15500.5  # J70 26001 -> 15500.5
$ws.Cells.Item(70, 12).Value = 15500.5  # L70 26001 -> 15500.5
$ws.Cells.Item(70, 14).Value = -16040.5  # N70 -26541 -> -16040.5
$ws.Cells.Item(73, 8).Value = 7730.6  # H73 8034 -> 7730.6
$ws.Cells.Item(73, 10).Value = 15500.5  # J73 26001 -> 15500.5
$ws.Cells.Item(73, 12).Value = 15500.5  # L73 26001 -> 15500.5
$ws.Cells.Item(73, 14).Value = -17372.5  # N73 -27873 -> -17372.5
$ws.Cells.Item(80, 8).Value = 3227.8  # H80 2848.611 -> 3227.8
$ws.Cells.Item(80, 9).Value = 3246.5  # I80 3007.7 -> 3246.5
$ws.Cells.Item(80, 10).Value = 3199.75  # J80 2649.75 -> 3199.75
$ws.Cells.Item(80, 11).Value = 3246.5  # K80 3007.7 -> 3246.5
$ws.Cells.Item(80, 12).Value = 3199.75  # L80 2649.75 -> 3199.75
$ws.Cells.Item(80, 13).Value = -2248.5  # M80 -2009.7 -> -2248.5
$ws.Cells.Item(80, 14).Value = -5195.75  # N80 -4645.75 -> -5195.75
$ws.Cells.Item(83, 8).Value = 3227.8  # H83 2848.611 -> 3227.8
$ws.Cells.Item(83, 9).Value = 3246.5  # I83 3007.7 -> 3246.5
$ws.Cells.Item(83, 10).Value = 3199.75  # J83 2649.75 -> 3199.75
$ws.Cells.Item(83, 11).Value = 16232.5  # K83 15038.5 -> 16232.5
$ws.Cells.Item(83, 12).Value = 15998.75  # L83 13248.75 -> 15998.75
$ws.Cells.Item(83, 13).Value = -11240.5  # M83 -10046.5 -> -11240.5
$ws.Cells.Item(83, 14).Value = -25982.75  # N83 -23232.75 -> -25982.75
$ws.Cells.Item(126, 8).Value = 16451.25  # H126 14757.056 -> 16451.25
$ws.Cells.Item(126, 9).Value = 19346.691  # I126 16058 -> 19346.691
$ws.Cells.Item(126, 10).Value = 3904.3333  # J126 4349.5 -> 3904.3333
$ws.Cells.Item(126, 11).Value = 58040.073  # K126 48174 -> 58040.073
$ws.Cells.Item(126, 12).Value = 11712.9999  # L126 13048.5 -> 11712.9999
$ws.Cells.Item(126, 13).Value = -55570.073  # M126 -45704 -> -55570.073
$ws.Cells.Item(126, 14).Value = -16652.9999  # N126 -17988.5 -> -16652.9999
$ws.Cells.Item(132, 8).Value = 3969.9312  # H132 4197.0835 -> 3969.9312
$ws.Cells.Item(132, 9).Value = 3161.2  # I132 3323.375 -> 3161.2
$ws.Cells.Item(132, 10).Value = 5767.1113  # J132 5944.5 -> 5767.1113
$ws.Cells.Item(132, 11).Value = 9483.599999999999  # K132 9970.125 -> 9483.599999999999
$ws.Cells.Item(132, 12).Value = 17301.3339  # L132 17833.5 -> 17301.3339
$ws.Cells.Item(132, 13).Value = -6953.599999999999  # M132 -7440.125 -> -6953.599999999999
$ws.Cells.Item(132, 14).Value = -22361.3339  # N132 -22893.5 -> -22361.3339
$ws.Cells.Item(134, 8).Value = 0  # H134 80000 -> 0
$ws.Cells.Item(134, 10).Value = 0  # J134 80000 -> 0
$ws.Cells.Item(134, 12).Value = 0  # L134 240000 -> 0
$ws.Cells.Item(134, 14).Value = $null  # clear N134 (was -245070)

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(55, 8).Value = 173.78261  # H55 181.4 -> 173.78261
$ws.Cells.Item(55, 9).Value = 159.85715  # I55 165.23077 -> 159.85715
$ws.Cells.Item(55, 10).Value = 195.44444  # J55 211.42857 -> 195.44444
$ws.Cells.Item(55, 11).Value = 159.85715  # K55 165.23077 -> 159.85715
$ws.Cells.Item(55, 12).Value = 195.44444  # L55 211.42857 -> 195.44444
$ws.Cells.Item(55, 13).Value = 13.14285000000001  # M55 7.769229999999993 -> 13.14285000000001
$ws.Cells.Item(55, 14).Value = -541.44444  # N55 -557.42857 -> -541.44444
$ws.Cells.Item(82, 8).Value = 1347.8572  # H82 1326.2142 -> 1347.8572
$ws.Cells.Item(82, 9).Value = 1188.1111  # I82 1135 -> 1188.1111
$ws.Cells.Item(82, 10).Value = 1635.4  # J82 1670.4 -> 1635.4
$ws.Cells.Item(82, 11).Value = 1188.1111  # K82 1135 -> 1188.1111
$ws.Cells.Item(82, 12).Value = 1635.4  # L82 1670.4 -> 1635.4
$ws.Cells.Item(82, 13).Value = -827.1111000000001  # M82 -774 -> -827.1111000000001
$ws.Cells.Item(82, 14).Value = -2357.4  # N82 -2392.4 -> -2357.4
$ws.Cells.Item(85, 8).Value = 1347.8572  # H85 1326.2142 -> 1347.8572
$ws.Cells.Item(85, 9).Value = 1188.1111  # I85 1135 -> 1188.1111
$ws.Cells.Item(85, 10).Value = 1635.4  # J85 1670.4 -> 1635.4
$ws.Cells.Item(85, 11).Value = 1188.1111  # K85 1135 -> 1188.1111
$ws.Cells.Item(85, 12).Value = 1635.4  # L85 1670.4 -> 1635.4
$ws.Cells.Item(85, 13).Value = 59.88889999999992  # M85 113 -> 59.88889999999992
$ws.Cells.Item(85, 14).Value = -4131.4  # N85 -4166.4 -> -4131.4
$ws.Cells.Item(122, 8).Value = 6943.6  # H122 8285.571 -> 6943.6
$ws.Cells.Item(122, 9).Value = 6722  # I122 7750 -> 6722
$ws.Cells.Item(122, 10).Value = 7276  # J122 8999.666999999999 -> 7276
$ws.Cells.Item(122, 11).Value = 20166  # K122 23250 -> 20166
$ws.Cells.Item(122, 12).Value = 21828  # L122 26999.001 -> 21828
$ws.Cells.Item(122, 13).Value = -17716  # M122 -20800 -> -17716
$ws.Cells.Item(122, 14).Value = -26728  # N122 -31899.001 -> -26728
$ws.Cells.Item(132, 8).Value = 3141.6897  # H132 3077.3547 -> 3141.6897
$ws.Cells.Item(132, 9).Value = 2311.8696  # I132 2386.52 -> 2311.8696
$ws.Cells.Item(132, 10).Value = 6322.6665  # J132 5955.8335 -> 6322.6665
$ws.Cells.Item(132, 11).Value = 6935.6088  # K132 7159.559999999999 -> 6935.6088
$ws.Cells.Item(132, 12).Value = 18967.9995  # L132 17867.5005 -> 18967.9995
$ws.Cells.Item(132, 13).Value = -4405.6088  # M132 -4629.559999999999 -> -4405.6088
$ws.Cells.Item(132, 14).Value = -24027.9995  # N132 -22927.5005 -> -24027.9995
$ws.Cells.Item(136, 8).Value = 6099.8887  # H136 5737.5 -> 6099.8887
$ws.Cells.Item(136, 9).Value = 5199.857  # I136 4566.6665 -> 5199.857
$ws.Cells.Item(136, 11).Value = 15599.571  # K136 13699.9995 -> 15599.571
$ws.Cells.Item(136, 13).Value = -13049.571  # M136 -11149.9995 -> -13049.571

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 6481.8623  # H81 6087.5483 -> 6481.8623
$ws.Cells.Item(81, 10).Value = 4355.7646  # J81 3936.2104 -> 4355.7646
$ws.Cells.Item(81, 12).Value = 8711.529200000001  # L81 7872.4208 -> 8711.529200000001
$ws.Cells.Item(81, 14).Value = -10833.5292  # N81 -9994.4208 -> -10833.5292
$ws.Cells.Item(84, 8).Value = 6481.8623  # H84 6087.5483 -> 6481.8623
$ws.Cells.Item(84, 10).Value = 4355.7646  # J84 3936.2104 -> 4355.7646
$ws.Cells.Item(84, 12).Value = 43557.64600000001  # L84 39362.104 -> 43557.64600000001
$ws.Cells.Item(84, 14).Value = -54165.64600000001  # N84 -49970.104 -> -54165.64600000001
$ws.Cells.Item(113, 8).Value = 2430  # H113 3003.75 -> 2430
$ws.Cells.Item(113, 9).Value = 1617.1177  # I113 2145.5833 -> 1617.1177
$ws.Cells.Item(113, 10).Value = 7036.3335  # J113 5578.25 -> 7036.3335
$ws.Cells.Item(113, 11).Value = 4851.3531  # K113 6436.749899999999 -> 4851.3531
$ws.Cells.Item(113, 12).Value = 21109.0005  # L113 16734.75 -> 21109.0005
$ws.Cells.Item(113, 13).Value = -2681.3531  # M113 -4266.749899999999 -> -2681.3531
$ws.Cells.Item(113, 14).Value = -25449.0005  # N113 -21074.75 -> -25449.0005
$ws.Cells.Item(122, 8).Value = 1979.421  # H122 2101.5625 -> 1979.421
$ws.Cells.Item(122, 9).Value = 1922.7222  # I122 2041.6666 -> 1922.7222
$ws.Cells.Item(122, 11).Value = 5768.1666  # K122 6124.9998 -> 5768.1666
$ws.Cells.Item(122, 13).Value = -3318.1666  # M122 -3674.9998 -> -3318.1666
$ws.Cells.Item(126, 8).Value = 1574.9  # H126 1691.3704 -> 1574.9
$ws.Cells.Item(126, 9).Value = 1385.96  # I126 1503.1364 -> 1385.96
$ws.Cells.Item(126, 11).Value = 4157.88  # K126 4509.4092 -> 4157.88
$ws.Cells.Item(126, 13).Value = -1687.88  # M126 -2039.4092 -> -1687.88
$ws.Cells.Item(132, 8).Value = 2072  # H132 2258.3948 -> 2072
$ws.Cells.Item(132, 9).Value = 1840.4849  # I132 1949.6 -> 1840.4849
$ws.Cells.Item(132, 10).Value = 3600  # J132 3416.375 -> 3600
$ws.Cells.Item(132, 11).Value = 5521.4547  # K132 5848.799999999999 -> 5521.4547
$ws.Cells.Item(132, 12).Value = 10800  # L132 10249.125 -> 10800
$ws.Cells.Item(132, 13).Value = -2991.4547  # M132 -3318.799999999999 -> -2991.4547
$ws.Cells.Item(132, 14).Value = -15860  # N132 -15309.125 -> -15860
$ws.Cells.Item(136, 8).Value = 1578.5883  # H136 1517.7894 -> 1578.5883
$ws.Cells.Item(136, 9).Value = 1367.5172  # I136 1350.3871 -> 1367.5172
$ws.Cells.Item(136, 10).Value = 2802.8  # J136 2259.1428 -> 2802.8
$ws.Cells.Item(136, 11).Value = 4102.5516  # K136 4051.1613 -> 4102.5516
$ws.Cells.Item(136, 12).Value = 8408.400000000001  # L136 6777.428400000001 -> 8408.400000000001
$ws.Cells.Item(136, 13).Value = -1552.5516  # M136 -1501.1613 -> -1552.5516
$ws.Cells.Item(136, 14).Value = -13508.4  # N136 -11877.4284 -> -13508.4

Write-Host "Applied 403 value updates and 7 cell clears across 8 sheets."